# Appointment Outcome Record List - add two new appointment outcome records
# (OR004 and OR005) as described in the commit: "edited doctor's record
# appointment outcome".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 5: OR004
$ws.Cells.Item(5, 1).Value = "OR004"
$ws.Cells.Item(5, 2).Value = "general"
$ws.Cells.Item(5, 3).Value = "test"
$ws.Cells.Item(5, 4).Value = "Pending"
$ws.Cells.Item(5, 5).Value = "Hello"

# New row 6: OR005
$ws.Cells.Item(6, 1).Value = "OR005"
$ws.Cells.Item(6, 2).Value = "test"
$ws.Cells.Item(6, 3).Value = "test"
$ws.Cells.Item(6, 4).Value = "Pending"
$ws.Cells.Item(6, 5).Value = "Hi test"

# Leave the selection on the row below the newly-entered data, matching
# where the cursor lands after typing the last row in Excel.
[void]$ws.Rows("6:6").Select()
